# Add a new "Swiss" test-data sheet, cloned from the existing "Czech" sheet,
# and populate it with the Switzerland market values.

$wb = $excel.ActiveWorkbook

# The "Czech" sheet (3rd tab) is the template for the new sheet.
$czechWs = $wb.Worksheets.Item(3)

# Before copying, change the Czech sheet's own selection to a "select all"
# state (mirrors what happened in the source workbook when the sheet was
# used as the basis for copying). Doing this while Czech is still the
# active sheet avoids forcing a tab switch.
[void]$czechWs.Range("A1:XFD1048576").Select()

# Duplicate the Czech sheet right after itself; Excel names the copy
# "Czech (2)" and makes it the active sheet/tab.
[void]$czechWs.Copy($null, $czechWs)
$newWs = $wb.Worksheets.Item(4)

# Rename the copy and fill in the Switzerland-specific values.
$newWs.Name = "Swiss"
$newWs.Range("B2").Value = "Switzerland Market"
$newWs.Range("B4").Value = "NGC-3476/T2344"

# Match the selection left on the new sheet.
[void]$newWs.Range("B2:B4").Select()
